$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Validation sheet: fill in accuracy/precision/recall/F-measure
# results for the two Naive Bayes / J48 rows, and move the
# selection (losing the tabSelected flag as Experimentation
# becomes the active sheet below).
# ---------------------------------------------------------------
$val = $wb.Worksheets.Item("Validation")
$val.Range("B2").Value = 0.801
$val.Range("C2").Value = 0.809
$val.Range("D2").Value = 0.801
$val.Range("E2").Value = 0.803
$val.Range("B3").Value = 0.782
$val.Range("C3").Value = 0.784
$val.Range("D3").Value = 0.782
$val.Range("E3").Value = 0.782

$val.Activate()
$val.Range("B4").Select()

# ---------------------------------------------------------------
# Experimentation sheet: add "Notes" and "Default" columns (C, D)
# describing the effect / default value of each hyper-parameter
# experiment, widen column C to match column B, and leave this
# sheet active/selected (C14) as the last-touched sheet.
# ---------------------------------------------------------------
$exp = $wb.Worksheets.Item("Experimentation")

$exp.Range("C1").Value = "Notes"
$exp.Range("C1").Font.Bold = $true
$exp.Range("D1").Value = "Default"
$exp.Range("D1").Font.Bold = $true

$exp.Range("B2").Value = "None"
$exp.Range("D2").Value = $false

$exp.Range("B3").Value = "None"
$exp.Range("D3").Value = $false

$exp.Range("B4").Value = "None"
$exp.Range("D4").Value = $false

$exp.Range("B5").Value = "Slightly different DT"
$exp.Range("D5").Value = $false

$exp.Range("B6").Value = "None"
$exp.Range("D6").Value = $true

$exp.Range("B7").Value = "Different DT, different metrics"
$exp.Range("C7").Value = "Lower CF = very slightly better metrics, shorter DT"
$exp.Range("D7").Value = 0.25

$exp.Range("B8").Value = "None"
$exp.Range("D8").Value = $false

$exp.Range("B9").Value = "None"
$exp.Range("D9").Value = $false

$exp.Range("B10").Value = "Significantly different DT, slightly different metrics"
$exp.Range("C10").Value = "After 4, metrics drop slowly but tree size drops rapidly"
$exp.Range("D10").Value = 2

$exp.Range("B11").Value = "Only applicable when reduced error pruning = true"
$exp.Range("C11").Value = "Only applicable when reduced error pruning = true"
$exp.Range("D11").Value = 3

$exp.Range("B12").Value = "Significantly different DT, worse metrics"
$exp.Range("C12").Value = "Shorter DT"
$exp.Range("D12").Value = $false

$exp.Range("B13").Value = "None"
$exp.Range("C13").Value = "Only applicable when reduced error pruning = true"
$exp.Range("D13").Value = 1

$exp.Range("B14").Value = "Larger tree, slightly worse metrics"
$exp.Range("D14").Value = $true

$exp.Range("B15").Value = "Larger tree, slightly worse metrics"
$exp.Range("C15").Value = "Removes confidence factor"
$exp.Range("D15").Value = $false

$exp.Range("B16").Value = "Same tree, same metrics, higher error rate"
$exp.Range("D16").Value = $false

$exp.Range("B17").Value = "None"
$exp.Range("D17").Value = $true

$exp.Columns.Item(3).ColumnWidth = 44.67

$exp.Activate()
$exp.Range("C14").Select()
